$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Registration Processor")
$ws.Activate()
$ws.Range("A21").Select()
Write-Host "Active sheet:" $wb.ActiveSheet.Name
